{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// The document contains a single 5-column table of two-digit multiplication\n// expressions (e.g. \"53\u00d774=\"). This script replaces each expression's text\n// with its updated value, in table (row-major) order, while leaving all\n// paragraph/run formatting (fonts, size, alignment, etc.) untouched because\n// only the `<w:t>` text content is rewritten.\n\nconst pairs = [\n  [\"53\u00d774=\", \"87\u00d720=\"],\n  [\"45\u00d761=\", \"64\u00d755=\"],\n  [\"20\u00d744=\", \"42\u00d755=\"],\n  [\"50\u00d711=\", \"38\u00d758=\"],\n  [\"69\u00d756=\", \"62\u00d741=\"],\n  [\"40\u00d727=\", \"40\u00d774=\"],\n  [\"80\u00d769=\", \"41\u00d744=\"],\n  [\"19\u00d785=\", \"94\u00d745=\"],\n  [\"87\u00d720=\", \"19\u00d721=\"],\n  [\"79\u00d731=\", \"87\u00d7100=\"],\n  [\"59\u00d777=\", \"29\u00d771=\"],\n  [\"53\u00d732=\", \"45\u00d738=\"],\n  [\"98\u00d775=\", \"12\u00d757=\"],\n  [\"96\u00d774=\", \"54\u00d722=\"],\n  [\"45\u00d748=\", \"57\u00d786=\"],\n  [\"39\u00d710=\", \"62\u00d769=\"],\n  [\"21\u00d761=\", \"95\u00d750=\"],\n  [\"66\u00d756=\", \"91\u00d794=\"],\n  [\"78\u00d712=\", \"79\u00d779=\"],\n  [\"42\u00d730=\", \"42\u00d754=\"],\n  [\"24\u00d715=\", \"29\u00d724=\"],\n  [\"63\u00d774=\", \"85\u00d757=\"],\n  [\"97\u00d796=\", \"65\u00d786=\"],\n  [\"53\u00d721=\", \"56\u00d711=\"],\n  [\"71\u00d784=\", \"54\u00d713=\"],\n  [\"21\u00d7100=\", \"24\u00d746=\"],\n  [\"81\u00d791=\", \"12\u00d725=\"],\n  [\"83\u00d712=\", \"35\u00d719=\"],\n  [\"83\u00d746=\", \"94\u00d786=\"],\n  [\"39\u00d729=\", \"98\u00d726=\"],\n  [\"72\u00d722=\", \"20\u00d799=\"],\n  [\"21\u00d765=\", \"26\u00d784=\"],\n  [\"42\u00d714=\", \"29\u00d759=\"],\n  [\"38\u00d787=\", \"96\u00d782=\"],\n  [\"25\u00d718=\", \"53\u00d782=\"],\n  [\"80\u00d759=\", \"38\u00d773=\"],\n  [\"87\u00d714=\", \"60\u00d710=\"],\n  [\"43\u00d769=\", \"27\u00d720=\"],\n  [\"51\u00d790=\", \"28\u00d729=\"],\n  [\"55\u00d716=\", \"42\u00d751=\"],\n  [\"81\u00d779=\", \"48\u00d744=\"],\n  [\"82\u00d754=\", \"91\u00d719=\"],\n  [\"36\u00d746=\", \"53\u00d749=\"],\n  [\"83\u00d781=\", \"52\u00d794=\"],\n  [\"45\u00d722=\", \"15\u00d746=\"],\n  [\"40\u00d712=\", \"58\u00d778=\"],\n  [\"89\u00d785=\", \"57\u00d746=\"],\n  [\"83\u00d766=\", \"44\u00d747=\"],\n  [\"32\u00d737=\", \"67\u00d787=\"],\n  [\"16\u00d738=\", \"11\u00d769=\"],\n  [\"53\u00d759=\", \"91\u00d711=\"],\n  [\"67\u00d784=\", \"21\u00d714=\"],\n  [\"100\u00d776=\", \"70\u00d737=\"],\n  [\"74\u00d797=\", \"53\u00d769=\"],\n  [\"31\u00d719=\", \"87\u00d730=\"],\n  [\"16\u00d783=\", \"10\u00d795=\"],\n  [\"22\u00d783=\", \"32\u00d735=\"],\n  [\"10\u00d789=\", \"100\u00d780=\"],\n  [\"28\u00d788=\", \"17\u00d762=\"],\n  [\"58\u00d765=\", \"42\u00d763=\"],\n  [\"50\u00d761=\", \"44\u00d767=\"],\n  [\"20\u00d738=\", \"70\u00d729=\"],\n  [\"78\u00d742=\", \"65\u00d787=\"],\n  [\"57\u00d754=\", \"50\u00d720=\"],\n  [\"29\u00d725=\", \"86\u00d718=\"],\n  [\"39\u00d777=\", \"88\u00d774=\"],\n  [\"23\u00d738=\", \"63\u00d758=\"],\n  [\"20\u00d777=\", \"29\u00d747=\"],\n  [\"72\u00d792=\", \"99\u00d718=\"],\n  [\"43\u00d795=\", \"89\u00d717=\"],\n  [\"58\u00d720=\", \"91\u00d751=\"],\n  [\"18\u00d722=\", \"42\u00d758=\"],\n  [\"29\u00d720=\", \"40\u00d731=\"],\n  [\"33\u00d777=\", \"42\u00d760=\"],\n  [\"95\u00d766=\", \"58\u00d779=\"],\n  [\"25\u00d795=\", \"84\u00d794=\"],\n  [\"64\u00d753=\", \"26\u00d770=\"],\n  [\"96\u00d788=\", \"59\u00d761=\"],\n  [\"72\u00d744=\", \"57\u00d751=\"],\n  [\"41\u00d765=\", \"82\u00d765=\"],\n  [\"52\u00d792=\", \"31\u00d778=\"],\n  [\"53\u00d720=\", \"36\u00d718=\"],\n  [\"25\u00d788=\", \"96\u00d730=\"],\n  [\"26\u00d713=\", \"33\u00d798=\"],\n  [\"22\u00d756=\", \"68\u00d784=\"],\n  [\"69\u00d763=\", \"84\u00d786=\"],\n  [\"19\u00d727=\", \"79\u00d728=\"],\n  [\"57\u00d795=\", \"37\u00d721=\"],\n  [\"42\u00d780=\", \"77\u00d792=\"],\n  [\"20\u00d789=\", \"22\u00d761=\"],\n  [\"45\u00d745=\", \"51\u00d793=\"],\n  [\"48\u00d742=\", \"59\u00d746=\"],\n  [\"29\u00d717=\", \"94\u00d795=\"],\n  [\"72\u00d783=\", \"94\u00d785=\"],\n  [\"59\u00d756=\", \"72\u00d712=\"],\n  [\"76\u00d787=\", \"32\u00d719=\"],\n  [\"12\u00d754=\", \"56\u00d780=\"],\n  [\"72\u00d733=\", \"47\u00d713=\"],\n  [\"25\u00d752=\", \"85\u00d720=\"],\n  [\"24\u00d785=\", \"33\u00d773=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected at least one table in the document body.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst columnCount = 5;\nif (pairs.length !== table.rowCount * columnCount) {\n  throw new Error(\n    `Expected ${table.rowCount * columnCount} cells but have ${pairs.length} replacement pairs.`\n  );\n}\n\n// Read every cell's current text first so we can confirm we are editing the\n// expected cell before writing (guards against an unexpected table shape).\nconst cells = [];\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < columnCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.load(\"value\");\n    cells.push(cell);\n  }\n}\nawait context.sync();\n\nfor (let i = 0; i < pairs.length; i++) {\n  const [oldText, newText] = pairs[i];\n  const cell = cells[i];\n  if (cell.value === oldText) {\n    cell.value = newText;\n  } else if (cell.value === newText) {\n    // Already updated - leave as is.\n  } else {\n    throw new Error(\n      `Cell ${i} text mismatch: expected \"${oldText}\" but found \"${cell.value}\".`\n    );\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument / $d is the open document.\n#\n# The document contains a single 5-column table of two-digit multiplication\n# expressions (e.g. \"53x74=\"). This script replaces each expression's text\n# with its updated value, in table (row-major) order, while leaving all\n# paragraph/run formatting (fonts, size, alignment, etc.) untouched because\n# only the cell Range.Text is rewritten (font/size live on the run and are\n# preserved by Word when Range.Text is reassigned).\n\n$pairs = @(\n  @(\"53\u00d774=\", \"87\u00d720=\"),\n  @(\"45\u00d761=\", \"64\u00d755=\"),\n  @(\"20\u00d744=\", \"42\u00d755=\"),\n  @(\"50\u00d711=\", \"38\u00d758=\"),\n  @(\"69\u00d756=\", \"62\u00d741=\"),\n  @(\"40\u00d727=\", \"40\u00d774=\"),\n  @(\"80\u00d769=\", \"41\u00d744=\"),\n  @(\"19\u00d785=\", \"94\u00d745=\"),\n  @(\"87\u00d720=\", \"19\u00d721=\"),\n  @(\"79\u00d731=\", \"87\u00d7100=\"),\n  @(\"59\u00d777=\", \"29\u00d771=\"),\n  @(\"53\u00d732=\", \"45\u00d738=\"),\n  @(\"98\u00d775=\", \"12\u00d757=\"),\n  @(\"96\u00d774=\", \"54\u00d722=\"),\n  @(\"45\u00d748=\", \"57\u00d786=\"),\n  @(\"39\u00d710=\", \"62\u00d769=\"),\n  @(\"21\u00d761=\", \"95\u00d750=\"),\n  @(\"66\u00d756=\", \"91\u00d794=\"),\n  @(\"78\u00d712=\", \"79\u00d779=\"),\n  @(\"42\u00d730=\", \"42\u00d754=\"),\n  @(\"24\u00d715=\", \"29\u00d724=\"),\n  @(\"63\u00d774=\", \"85\u00d757=\"),\n  @(\"97\u00d796=\", \"65\u00d786=\"),\n  @(\"53\u00d721=\", \"56\u00d711=\"),\n  @(\"71\u00d784=\", \"54\u00d713=\"),\n  @(\"21\u00d7100=\", \"24\u00d746=\"),\n  @(\"81\u00d791=\", \"12\u00d725=\"),\n  @(\"83\u00d712=\", \"35\u00d719=\"),\n  @(\"83\u00d746=\", \"94\u00d786=\"),\n  @(\"39\u00d729=\", \"98\u00d726=\"),\n  @(\"72\u00d722=\", \"20\u00d799=\"),\n  @(\"21\u00d765=\", \"26\u00d784=\"),\n  @(\"42\u00d714=\", \"29\u00d759=\"),\n  @(\"38\u00d787=\", \"96\u00d782=\"),\n  @(\"25\u00d718=\", \"53\u00d782=\"),\n  @(\"80\u00d759=\", \"38\u00d773=\"),\n  @(\"87\u00d714=\", \"60\u00d710=\"),\n  @(\"43\u00d769=\", \"27\u00d720=\"),\n  @(\"51\u00d790=\", \"28\u00d729=\"),\n  @(\"55\u00d716=\", \"42\u00d751=\"),\n  @(\"81\u00d779=\", \"48\u00d744=\"),\n  @(\"82\u00d754=\", \"91\u00d719=\"),\n  @(\"36\u00d746=\", \"53\u00d749=\"),\n  @(\"83\u00d781=\", \"52\u00d794=\"),\n  @(\"45\u00d722=\", \"15\u00d746=\"),\n  @(\"40\u00d712=\", \"58\u00d778=\"),\n  @(\"89\u00d785=\", \"57\u00d746=\"),\n  @(\"83\u00d766=\", \"44\u00d747=\"),\n  @(\"32\u00d737=\", \"67\u00d787=\"),\n  @(\"16\u00d738=\", \"11\u00d769=\"),\n  @(\"53\u00d759=\", \"91\u00d711=\"),\n  @(\"67\u00d784=\", \"21\u00d714=\"),\n  @(\"100\u00d776=\", \"70\u00d737=\"),\n  @(\"74\u00d797=\", \"53\u00d769=\"),\n  @(\"31\u00d719=\", \"87\u00d730=\"),\n  @(\"16\u00d783=\", \"10\u00d795=\"),\n  @(\"22\u00d783=\", \"32\u00d735=\"),\n  @(\"10\u00d789=\", \"100\u00d780=\"),\n  @(\"28\u00d788=\", \"17\u00d762=\"),\n  @(\"58\u00d765=\", \"42\u00d763=\"),\n  @(\"50\u00d761=\", \"44\u00d767=\"),\n  @(\"20\u00d738=\", \"70\u00d729=\"),\n  @(\"78\u00d742=\", \"65\u00d787=\"),\n  @(\"57\u00d754=\", \"50\u00d720=\"),\n  @(\"29\u00d725=\", \"86\u00d718=\"),\n  @(\"39\u00d777=\", \"88\u00d774=\"),\n  @(\"23\u00d738=\", \"63\u00d758=\"),\n  @(\"20\u00d777=\", \"29\u00d747=\"),\n  @(\"72\u00d792=\", \"99\u00d718=\"),\n  @(\"43\u00d795=\", \"89\u00d717=\"),\n  @(\"58\u00d720=\", \"91\u00d751=\"),\n  @(\"18\u00d722=\", \"42\u00d758=\"),\n  @(\"29\u00d720=\", \"40\u00d731=\"),\n  @(\"33\u00d777=\", \"42\u00d760=\"),\n  @(\"95\u00d766=\", \"58\u00d779=\"),\n  @(\"25\u00d795=\", \"84\u00d794=\"),\n  @(\"64\u00d753=\", \"26\u00d770=\"),\n  @(\"96\u00d788=\", \"59\u00d761=\"),\n  @(\"72\u00d744=\", \"57\u00d751=\"),\n  @(\"41\u00d765=\", \"82\u00d765=\"),\n  @(\"52\u00d792=\", \"31\u00d778=\"),\n  @(\"53\u00d720=\", \"36\u00d718=\"),\n  @(\"25\u00d788=\", \"96\u00d730=\"),\n  @(\"26\u00d713=\", \"33\u00d798=\"),\n  @(\"22\u00d756=\", \"68\u00d784=\"),\n  @(\"69\u00d763=\", \"84\u00d786=\"),\n  @(\"19\u00d727=\", \"79\u00d728=\"),\n  @(\"57\u00d795=\", \"37\u00d721=\"),\n  @(\"42\u00d780=\", \"77\u00d792=\"),\n  @(\"20\u00d789=\", \"22\u00d761=\"),\n  @(\"45\u00d745=\", \"51\u00d793=\"),\n  @(\"48\u00d742=\", \"59\u00d746=\"),\n  @(\"29\u00d717=\", \"94\u00d795=\"),\n  @(\"72\u00d783=\", \"94\u00d785=\"),\n  @(\"59\u00d756=\", \"72\u00d712=\"),\n  @(\"76\u00d787=\", \"32\u00d719=\"),\n  @(\"12\u00d754=\", \"56\u00d780=\"),\n  @(\"72\u00d733=\", \"47\u00d713=\"),\n  @(\"25\u00d752=\", \"85\u00d720=\"),\n  @(\"24\u00d785=\", \"33\u00d773=\"),\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$columnCount = 5\n$rowCount = $t.Rows.Count\nif ($pairs.Length -ne ($rowCount * $columnCount)) {\n  throw \"Expected $($rowCount * $columnCount) cells but have $($pairs.Length) replacement pairs.\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $columnCount; $c++) {\n    $cell = $t.Cell($r, $c)\n    $rng = $cell.Range\n    $current = $rng.Text.TrimEnd([char]13, [char]7)\n    $oldText = $pairs[$i][0]\n    $newText = $pairs[$i][1]\n    if ($current -eq $oldText) {\n      $rng.Text = $newText\n    } elseif ($current -eq $newText) {\n      # Already updated - leave as is.\n    } else {\n      throw \"Cell ($r,$c) text mismatch: expected '$oldText' but found '$current'.\"\n    }\n    $i++\n  }\n}\n"}
